# Edit script: applies the WS2-2 Binary Representation title/marks changes.
#
# 1) Split the "WORKSHEET 2-2" run (first paragraph) into two runs:
#       "WORKSHEET 2-" + "2"
# 2) Add a new run "BINARY REPRESENTATION" to the (previously empty)
#    second, right-aligned title paragraph, keeping its pPr unchanged.
# 3) Turn the empty paragraph that follows "Please do the work ..." into
#    a bold/italic "MAX MARKS: 21" line (both pPr/rPr and run rPr change
#    from bCs-only to b/i/iCs).

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="2048"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Step 1: split "WORKSHEET 2-2" into "WORKSHEET 2-" + "2" (same rPr).
# Locate it with Find so the preceding drawing/run in the paragraph is
# left completely untouched.
# ---------------------------------------------------------------------
$titleRng = $d.Content
$found = $titleRng.Find.Execute("WORKSHEET 2-2", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'WORKSHEET 2-2' text"
}

$titleBody = '<w:p><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiBold" w:hAnsi="Bahnschrift SemiBold"/><w:b/><w:bCs/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr><w:t>WORKSHEET 2-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiBold" w:hAnsi="Bahnschrift SemiBold"/><w:b/><w:bCs/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr><w:t>2</w:t></w:r></w:p>'
[void]$titleRng.InsertXML($pkgHeader + $titleBody + $pkgFooter)

# ---------------------------------------------------------------------
# Step 2: the paragraph right after the title (currently empty,
# right-aligned) gets a "BINARY REPRESENTATION" run. Its pPr (spacing,
# alignment, rPr) is reproduced unchanged.
# ---------------------------------------------------------------------
$subtitlePara = $d.Paragraphs.Item(2)
$subtitleBody = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Bahnschrift SemiBold" w:hAnsi="Bahnschrift SemiBold"/><w:b/><w:bCs/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift SemiBold" w:hAnsi="Bahnschrift SemiBold"/><w:b/><w:bCs/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr><w:t>BINARY REPRESENTATION</w:t></w:r></w:p>'
[void]$subtitlePara.Range.InsertXML($pkgHeader + $subtitleBody + $pkgFooter)

# ---------------------------------------------------------------------
# Step 3: find the empty paragraph right after "Please do the work..."
# and give it a bold/italic "MAX MARKS: 21" run; pPr rPr switches from
# bCs-only to b/i/iCs.
# ---------------------------------------------------------------------
$marksRng = $d.Content
$found2 = $marksRng.Find.Execute("Please do the work either on the sheet or by hand and show it, that way I know you", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'Please do the work...' paragraph"
}
$precedingIndex = $marksRng.Paragraphs.Item(1).Index
$marksPara = $d.Paragraphs.Item($precedingIndex + 1)
$marksBody = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Bahnschrift" w:hAnsi="Bahnschrift"/><w:b/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Bahnschrift" w:hAnsi="Bahnschrift"/><w:b/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>MAX MARKS: 21</w:t></w:r></w:p>'
[void]$marksPara.Range.InsertXML($pkgHeader + $marksBody + $pkgFooter)

Write-Output "Edits applied successfully."
